# Update computed results for Case_0_133 (380 kV case) - pl_mw sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "C" = 0.0481268431561972; "D" = 0.1265004971426862; "E" = 0.1492802303412297; "F" = 2.035767034657027; "G" = 1.371116451407573; "H" = 1.260690378233804; "J" = 0.2018840846851759; "K" = 2.193010561205085 }
  3 = @{ "C" = 0.04277889447789107; "D" = 0.122595842664694; "E" = 0.1457044048105054; "F" = 2.036353246782852; "G" = 1.371763944753155; "H" = 1.269214961302609; "J" = 0.1977992357084304; "K" = 1.994407359270156 }
  4 = @{ "C" = 0.0395133323356589; "D" = 0.1202207000867048; "E" = 0.1435785126998539; "F" = 2.038324911203063; "G" = 1.373479708666764; "H" = 1.275346907733805; "J" = 0.1954213259261763; "K" = 1.872871147825833 }
  5 = @{ "C" = 0.03818704151484553; "D" = 0.1192584953131686; "E" = 0.142729703060926; "F" = 2.039531971766138; "G" = 1.374508471525118; "H" = 1.278070766466456; "J" = 0.1944849192315203; "K" = 1.823447069488452 }
  6 = @{ "C" = 0.03796707796809073; "D" = 0.1190990674023098; "E" = 0.142589816011462; "F" = 2.039756729837777; "G" = 1.374699146914068; "H" = 1.278536631945443; "J" = 0.1943313954653689; "K" = 1.815246473531488 }
  7 = @{ "C" = 0.03949542761668567; "D" = 0.1202077003427462; "E" = 0.1435669944791549; "F" = 2.038339558174002; "G" = 1.373492251009736; "H" = 1.275382732436483; "J" = 0.1954085653474564; "K" = 1.872204179896187 }
  8 = @{ "C" = 0.04627905226931261; "D" = 0.125149577427706; "E" = 0.1480328142067115; "F" = 2.035633620919626; "G" = 1.37106496740941; "H" = 1.263442875174007; "J" = 0.2004485244239902; "K" = 2.124447977172395 }
  9 = @{ "C" = 0.0597305038016458; "D" = 0.1350153293305851; "E" = 0.1573444153797894; "F" = 2.043194986767801; "G" = 1.376852066561383; "H" = 1.247186069153173; "J" = 0.2113716764561673; "K" = 2.622330425073983 }
  10 = @{ "C" = 0.06971215144390897; "D" = 0.1423678687233547; "E" = 0.1645262289568166; "F" = 2.056709202240356; "G" = 1.387657998659535; "H" = 1.239652934691776; "J" = 0.2200412888365548; "K" = 2.99014514886187 }
  11 = @{ "C" = 0.07427636467801335; "D" = 0.1457349384901221; "E" = 0.1678680589961985; "F" = 2.064610466183808; "G" = 1.394024192947199; "H" = 1.237194075507091; "J" = 0.2241276029385517; "K" = 3.157926876301474 }
  12 = @{ "C" = 0.07600822352368652; "D" = 0.1470131266728743; "E" = 0.1691443161185546; "F" = 2.067856615248459; "G" = 1.396645655487447; "H" = 1.236403010551101; "J" = 0.225695641785677; "K" = 3.22152816832147 }
  13 = @{ "C" = 0.07563507974300876; "D" = 0.1467377067889686; "E" = 0.1688689715240201; "F" = 2.067146166064333; "G" = 1.396071670451278; "H" = 1.236567138025237; "J" = 0.2253570160956855; "K" = 3.207827560848045 }
  14 = @{ "C" = 0.07441877516293971; "D" = 0.1458400329896961; "E" = 0.1679728412724373; "F" = 2.064872424821843; "G" = 1.394235626745058; "H" = 1.237126182356604; "J" = 0.2242561918818637; "K" = 3.163158071661258 }
  15 = @{ "C" = 0.07367421130736318; "D" = 0.1452905904852457; "E" = 0.1674253400650088; "F" = 2.063512843635706; "G" = 1.393138503030144; "H" = 1.237486877481814; "J" = 0.2235845978727582; "K" = 3.135805306680822 }
  16 = @{ "C" = 0.06941435180873157; "D" = 0.1421482672906649; "E" = 0.1643093394491686; "F" = 2.056228288850036; "G" = 1.387271323704994; "H" = 1.239833178445394; "J" = 0.2197771184672206; "K" = 2.979189467282708 }
  17 = @{ "C" = 0.06680716526543051; "D" = 0.140226235373035; "E" = 0.1624169454905768; "F" = 2.052209966421401; "G" = 1.384045030950432; "H" = 1.241521092037857; "J" = 0.2174779505742066; "K" = 2.883228552633284 }
  18 = @{ "C" = 0.06530978692138945; "D" = 0.1391228420182813; "E" = 0.1613355282799844; "F" = 2.050063671457835; "G" = 1.382325753997605; "H" = 1.242582996708421; "J" = 0.2161689300241818; "K" = 2.82807777318817 }
  19 = @{ "C" = 0.06480317598401086; "D" = 0.1387496161135005; "E" = 0.1609705870211684; "F" = 2.049365241918835; "G" = 1.3817669984156; "H" = 1.242958153341306; "J" = 0.2157280149628065; "K" = 2.809412131242482 }
  20 = @{ "C" = 0.06708447557952013; "D" = 0.140430621189779; "E" = 0.1626176655064384; "F" = 2.052620639245475; "G" = 1.384374342885224; "H" = 1.241331980116854; "J" = 0.2177213129604638; "K" = 2.893439266377811 }
  21 = @{ "C" = 0.07477593769561963; "D" = 0.146103616499957; "E" = 0.168235763718549; "F" = 2.065533365605688; "G" = 1.394769181209796; "H" = 1.236958169678985; "J" = 0.2245789692225486; "K" = 3.176276794336843 }
  22 = @{ "C" = 0.07982314908252874; "D" = 0.149829590642625; "E" = 0.1719703542664988; "F" = 2.07545465292516; "G" = 1.402792065926093; "H" = 1.234916361510329; "J" = 0.2291812253254477; "K" = 3.361512049596115 }
  23 = @{ "C" = 0.07712745458525205; "D" = 0.1478393102069901; "E" = 0.1699713747586316; "F" = 2.070023217189771; "G" = 1.398396913823717; "H" = 1.235931097991113; "J" = 0.2267138481156081; "K" = 3.26261339435149 }
  24 = @{ "C" = 0.06695909883741535; "D" = 0.1403382133482296; "E" = 0.1625268996027671; "F" = 2.052434463696073; "G" = 1.384225038857949; "H" = 1.241417192657167; "J" = 0.2176112489639763; "K" = 2.888822945279458 }
  25 = @{ "C" = 0.0560747625828526; "D" = 0.1323279258519818; "E" = 0.154765776909116; "F" = 2.039760085166591; "G" = 1.374144553698301; "H" = 1.250812848154112; "J" = 0.2083042811477895; "K" = 2.487289085474686 }
}

foreach ($row in $data.Keys) {
  $rowData = $data[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = $rowData[$col]
  }
}

Write-Output "Updated $($data.Count) rows across columns C,D,E,F,G,H,J,K"
